# Update the cryptocurrency price/volume table with the latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '42.762.18'
$ws.Range("E2").Value = '  -0.07%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.303.11'
$ws.Range("E3").Value = '  +0.43%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5: BNB
$ws.Range("D5").Value = '301.32'
$ws.Range("E5").Value = '  -1.46%  '

# Row 6: Solana
$ws.Range("D6").Value = '96.10'
$ws.Range("E6").Value = '  -0.46%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +0.12%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.07%  '

# Row 9: Cardano
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  -1.30%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '34.72'
$ws.Range("E10").Value = '  -2.38%  '

# Row 11: Chainlink
$ws.Range("D11").Value = '19.26'
$ws.Range("E11").Value = '  +4.99%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = '0.0789'
$ws.Range("E12").Value = '  -0.12%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -0.40%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '6.77'
$ws.Range("E14").Value = '  +0.63%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.658.21'
$ws.Range("E15").Value = '  +0.24%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.314.20'
$ws.Range("E16").Value = '  +0.86%  '

# Row 17: Polygon
$ws.Range("D17").Value = '0.786'
$ws.Range("E17").Value = '  +0.91%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '42.659.01'
$ws.Range("E18").Value = '  -0.14%  '

# Row 19: InternetComputer(DFINITY)
$ws.Range("D19").Value = '12.28'
$ws.Range("E19").Value = '  -5.76%  '

# Row 20: ShibaInu
$ws.Range("D20").Value = '0.0₃0892'
$ws.Range("E20").Value = '  -0.64%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '6.03'
$ws.Range("E21").Value = '  +0.58%  '

# Row 22: Litecoin
$ws.Range("D22").Value = '67.78'
$ws.Range("E22").Value = '  +0.90%  '

# Row 23: ImmutableX
$ws.Range("E23").Value = '  +6.92%  '

# Row 24: BitcoinCash
$ws.Range("D24").Value = '235.22'
$ws.Range("E24").Value = '  -0.35%  '

# Row 25: Dai
$ws.Range("E25").Value = '  +0.16%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  -2.13%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '24.40'
$ws.Range("E27").Value = '  -2.93%  '

# Row 28: Toncoin
$ws.Range("E28").Value = '  +14.79%  '

# Row 29: Monero
$ws.Range("D29").Value = '165.05'
$ws.Range("E29").Value = '  -0.96%  '

# Row 30: Cosmos
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  +0.31%  '

# Row 31: InjectiveProtocol
$ws.Range("D31").Value = '32.23'
$ws.Range("E31").Value = '  -2.99%  '

# Row 32: FirstDigitalUSD
$ws.Range("E32").Value = '  -0.09%  '

# Row 33: Filecoin
$ws.Range("D33").Value = '4.98'
$ws.Range("E33").Value = '  +0.25%  '

# Row 34: Celestia
$ws.Range("D34").Value = '17.59'
$ws.Range("E34").Value = '  -0.53%  '

# Row 35: RenderToken
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -6.08%  '

# Row 36: Hedera
$ws.Range("D36").Value = '0.0705'
$ws.Range("E36").Value = '  +2.18%  '

# Row 37: WEMIXToken
$ws.Range("E37").Value = '  -3.00%  '

# Row 38: Kaspa
$ws.Range("E38").Value = '  -0.55%  '

# Row 39: ARBITRUM
$ws.Range("E39").Value = '  +0.43%  '

# Row 40: LidoDAOToken
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  +0.52%  '

# Row 41: Stellar
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.108'
$ws.Range("E41").Value = '  -0.73%  '

# Row 42: EnergySwap
$ws.Range("D42").Value = '20.39'
$ws.Range("E42").Value = '  +12.38%  '

# Row 43: Maker
$ws.Range("D43").Value = '1.972.87'
$ws.Range("E43").Value = '  -1.32%  '

# Row 44: FraxShare
$ws.Range("D44").Value = '10.45'
$ws.Range("E44").Value = '  +4.83%  '

# Row 45: VeChain
$ws.Range("E45").Value = '  -0.14%  '

# Row 46: ApeXProtocol
$ws.Range("D46").Value = '2.02'
$ws.Range("E46").Value = '  -2.68%  '

# Row 47: NEARProtocol
$ws.Range("D47").Value = '2.78'
$ws.Range("E47").Value = '  +0.39%  '

# Row 48: RocketPoolETH
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.525.50'
$ws.Range("E48").Value = '  +0.18%  '

# Row 49: MultiversX
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = '53.32'
$ws.Range("E49").Value = '  -0.68%  '

# Row 50: HuobiToken
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '2.80'
$ws.Range("E50").Value = '  -1.18%  '

# Row 51: BitcoinSV
$ws.Range("D51").Value = '71.54'
$ws.Range("E51").Value = '  +0.37%  '

